$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right after the current "Spring_2016" row (which is
# row 2) so the new rows inherit that row's formatting (bold custom row
# style), rather than the plain header-row formatting they'd inherit if
# inserted immediately above row 2.
$ws.Range("A3:F4").EntireRow.Insert()

# Old row 2 data (Spring_2016 / current / priority 0.6) now slides down
# conceptually: row 2 keeps its place but becomes the new "Spring_2017"
# entry, row 3 becomes "Fall_2016" / current, and row 4 becomes the old
# Spring_2016 data with an added lastmod date and a lower priority.

# Row 2: new upcoming semester "Spring_2017" / status "future"
$ws.Range("A2").Value = "Fall_2016"
$ws.Range("A2").Value = "Spring_2017"
$ws.Range("B2").Value = "future"
$ws.Range("C2").Value = $null
$ws.Range("D2").Value = $null

# Row 3: "Fall_2016" / status "current" (priority 0.6)
$ws.Range("A3").Value = "Fall_2016"
$ws.Range("B3").Value = "current"
$ws.Range("D3").Value = 0.6

# Row 4: old "Spring_2016" row, now "past" with a real lastmod date and a
# lower priority than before.
$ws.Range("A4").Value = "Spring_2016"
$ws.Range("B4").Value = "past"
$ws.Range("C4").Value = 42545
$ws.Range("D4").Value = 0.3

# Row 5 (Fall_2015) and row 6 (Spring_2015) priorities step down as they
# age further away from "current".
$ws.Range("D5").Value = 0.1
$ws.Range("D6").Value = 0.05

$ws.Range("A25").Select()
